$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-shuffles the data rows 2-9 (only columns D, M, N, O, P, Q, S vary
# between rows; all other columns are identical across the block). Capture the
# "before" values for those columns per row, then write them back out under
# the new row mapping.

$cols = @("D", "M", "N", "O", "P", "Q", "S")

$before = @{}
for ($r = 2; $r -le 9; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: new row -> source (old) row whose D/M/N/O/P/Q/S values it receives
$mapping = @{
    2 = 4
    3 = 8
    4 = 9
    5 = 3
    6 = 2
    7 = 7
    8 = 6
    9 = 5
}

foreach ($newRow in $mapping.Keys) {
    $srcRow = $mapping[$newRow]
    $src = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
